# "prepare for new R version"
#
# The "meta" sheet stores key/value configuration pairs for the chart
# (column A = key, column B = value). This change:
#   1. Appends a "(%)" suffix to the y_title value.
#   2. Introduces a new "box_median_lab_suffix" = "\s%" key/value pair,
#      inserted right after "box_median_col" / "rose" (i.e. as the new
#      row 12), pushing the rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# 1. Update y_title's value to mention it is now a percentage.
$ws.Range("B6").Value = "toename in de kans om geraakt te worden t.o.v. het gemiddelde (%)"

# 2. Insert a new row after "box_median_col" (row 11) / before "hline_dash"
#    (old row 12), and populate it.
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "box_median_lab_suffix"
$ws.Range("B12").Value = "\s%"
